# Update the acquisition timestamp (column A) for all data rows on the
# "ランサーズ" sheet from 2025-11-07 06:27:17 to 2025-11-07 06:35:04.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-07 06:35:04"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
